# Formed the consolidated report
# Update the "Absent" (column H) values so that they correctly reflect
# Absent = 1 - Real (column E) for each attendance row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
